$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 currently only has formulas in E20/F20 (carried down from the shared
# range). Pull the cell formatting from row 19 (the last fully-filled entry)
# so the new row matches styles used by the rest of the log (D20 in
# particular needs the time-formatted style, not the column default).
$ws.Range("B19:H19").Copy() | Out-Null
$ws.Range("B20:H20").PasteSpecial(-4122) | Out-Null

# New timesheet entry for 45211 (2023-10-12), 5:30 PM - 7:15 PM
$ws.Range("B20").Value = 45211
$ws.Range("C20").Value = 0.72916666666666663
$ws.Range("D20").Value = 0.80208333333333337

# Journal notes for the new session
$ws.Range("G20").Value = "I got the web containers mounting a react template that is stored as json instead of using npx create-react-app. I don't know if it's faster, but it is more controllable."
$ws.Range("H20").Value = "Next task is the file heirarchy viewer/selector. First view, then select opens it on the code editor"

# Move the active selection to where editing left off
$ws.Range("H21").Select() | Out-Null
